$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original values of the columns that change (D, J, K, L, M, O, P)
# for rows 2-4 before overwriting anything, since the new values are a
# cyclic shift of the old ones (row2 -> row4, row3 -> row2, row4 -> row3).
$cols = @("D", "J", "K", "L", "M", "O", "P")

$orig = @{}
foreach ($row in 2..4) {
    $orig[$row] = @{}
    foreach ($col in $cols) {
        $orig[$row][$col] = $ws.Range("$col$row").Value2
    }
}

# New row 2 = old row 3, new row 3 = old row 4, new row 4 = old row 2
$mapping = @{ 2 = 3; 3 = 4; 4 = 2 }

foreach ($destRow in 2..4) {
    $srcRow = $mapping[$destRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $orig[$srcRow][$col]
    }
}
